$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The old B column (the "4"/StimSet repeat count) is being dropped; the
# sheet becomes a single column of image-file stimuli again. Clear the old
# A1:B4 block and rebuild it as A1:A6.
$ws.Range("A1:B4").Clear()

$ws.Range("A1").Value = "ImageFile"
$ws.Range("A2").Value = "Stimuli/325.jpg"
$ws.Range("A3").Value = "Stimuli/1300.jpg"
$ws.Range("A4").Value = "Stimuli/2457.jpg"
$ws.Range("A5").Value = "Stimuli/2683.jpg"
$ws.Range("A6").Value = "Stimuli/6314.jpg"

$ws.Range("A7").Select()
